# Apply cryptos list update (values scraped on Sun Apr 16 15:22:15 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($Sheet, $Address, $Text)
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = '@'
    $cell.Value = $Text
    $cell.Style = 'Normal'
}

$ws.Range('D2').Value = '30.659.30'
$ws.Range('D3').Value = '2.112.85'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +1.10%  '
Set-TextCellValue $ws 'D5' '350.00'
$ws.Range('E5').Value = '  +4.18%  '
Set-TextCellValue $ws 'D6' '1.011'
$ws.Range('E6').Value = '  +1.00%  '
Set-TextCellValue $ws 'D7' '0.5262'
$ws.Range('E7').Value = '  +0.27%  '
Set-TextCellValue $ws 'D8' '0.4509'
$ws.Range('E8').Value = '  -1.99%  '
Set-TextCellValue $ws 'D9' '53.75'
$ws.Range('E9').Value = '  +0.84%  '
Set-TextCellValue $ws 'D10' '0.09008'
$ws.Range('E10').Value = '  +0.55%  '
$ws.Range('E11').Value = '  -0.75%  '
Set-TextCellValue $ws 'D12' '24.47'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '2.125.98'
$ws.Range('E13').Value = '  +1.01%  '
Set-TextCellValue $ws 'D14' '6.811'
$ws.Range('E14').Value = '  -0.01%  '
Set-TextCellValue $ws 'D15' '8.015'
$ws.Range('E15').Value = '  +0.51%  '
Set-TextCellValue $ws 'D16' '99.88'
$ws.Range('E16').Value = '  +3.36%  '
Set-TextCellValue $ws 'D17' '0.00001173'
$ws.Range('E17').Value = '  +3.44%  '
Set-TextCellValue $ws 'D20' '19.34'
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('E21').Value = '  +1.00%  '
Set-TextCellValue $ws 'D22' '6.286'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').Value = '30.716.47'
$ws.Range('E23').Value = '  +0.55%  '
Set-TextCellValue $ws 'D24' '12.84'
$ws.Range('E24').Value = '  +3.82%  '
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').Value = '2.369.49'
$ws.Range('E26').Value = '  +1.59%  '
Set-TextCellValue $ws 'D27' '22.29'
$ws.Range('E27').Value = '  -0.30%  '
Set-TextCellValue $ws 'D28' '165.43'
$ws.Range('E28').Value = '  +0.89%  '
$ws.Range('E29').Value = '  -1.85%  '
Set-TextCellValue $ws 'D30' '135.54'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('E31').Value = '  -2.63%  '
$ws.Range('E32').Value = '  -0.09%  '
Set-TextCellValue $ws 'D33' '1.640'
$ws.Range('E33').Value = '  -3.84%  '
Set-TextCellValue $ws 'D34' '6.335'
$ws.Range('E34').Value = '  +1.98%  '
Set-TextCellValue $ws 'D35' '4.020'
$ws.Range('E35').Value = '  +2.37%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCellValue $ws 'D36' '10.21'
$ws.Range('E36').Value = '  -2.50%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCellValue $ws 'D37' '5.895'
$ws.Range('E37').Value = '  +5.99%  '
Set-TextCellValue $ws 'D38' '0.02650'
$ws.Range('E38').Value = '  +2.84%  '
Set-TextCellValue $ws 'D39' '0.06830'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  +0.37%  '
Set-TextCellValue $ws 'D41' '12.53'
$ws.Range('E41').Value = '  -2.81%  '
Set-TextCellValue $ws 'D42' '0.6881'
$ws.Range('E42').Value = '  -0.26%  '
Set-TextCellValue $ws 'D43' '1.276'
$ws.Range('E43').Value = '  +2.19%  '
Set-TextCellValue $ws 'D44' '14.74'
$ws.Range('E44').Value = '  +5.24%  '
Set-TextCellValue $ws 'D45' '2.319'
$ws.Range('E45').Value = '  -1.35%  '
Set-TextCellValue $ws 'D46' '0.6426'
$ws.Range('E46').Value = '  +0.46%  '
Set-TextCellValue $ws 'D47' '3.763'
$ws.Range('E47').Value = '  +2.56%  '
$ws.Range('E48').Value = '  -0.68%  '
Set-TextCellValue $ws 'D49' '1.248'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('E50').Value = '  +2.17%  '
$ws.Range('B51').Value = 'WEMIXTOKEN'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCellValue $ws 'D51' '1.193'
$ws.Range('E51').Value = '  -0.21%  '
